# "Add new icons and fix API endpoint"
#
# The quiz export sheet used to lead with an "id_quiz" column (A) followed
# by question_text / answer_text1 / is_correct1 / answer_text2 /
# is_correct2 in B:F. The API endpoint feeding this sheet no longer
# returns id_quiz, so that column is dropped: every other column's data
# shifts one slot to the left (B:F -> A:E) and the now-unused trailing
# column F is cleared. The freed-up column A (now holding the question
# text) is given a wider custom width to fit its new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header row B1:F1 into A1:E1 (drops the old "id_quiz" header in A1),
# then blank out the now-trailing F1.
$ws.Range("B1:F1").Copy($ws.Range("A1"))
$ws.Range("F1").Clear()

# Shift the data row B2:F2 into A2:E2 (drops the old id_quiz value in A2,
# carries the question-text cell's style along with it), then blank out
# the now-trailing F2.
$ws.Range("B2:F2").Copy($ws.Range("A2"))
$ws.Range("F2").Clear()

# Give the new first column (question text) a wider custom width.
$ws.Columns.Item(1).ColumnWidth = 20.1666666666666667

# Restore the active selection.
$ws.Range("C9").Select()
